# Update tracklist workbook:
#  - rename the two "kendricklamar4" defined names to "mkiwanuka2"
#  - replace the Kendrick Lamar track data on Sheet1 / Sheet3 with the
#    Michael Kiwanuka track data (incl. the two new rows 10 & 11)
#  - resize the data columns on Sheet1 / Sheet3 to their new widths
#  - move the Sheet2 selection down to the new last data row (K14)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename defined names (kendricklamar4 -> mkiwanuka2) on every sheet
#    that has one. Renaming shifts the collection order (names are kept
#    sorted), so re-scan from the top until nothing starting with the
#    old name is left.
# ---------------------------------------------------------------------
$stillOld = $true
while ($stillOld) {
    $stillOld = $false
    for ($i = 1; $i -le $wb.Names.Count; $i++) {
        $nm = $wb.Names.Item($i)
        if ($nm.Name -like "*kendricklamar4*") {
            $nm.Name = "mkiwanuka2"
            $stillOld = $true
            break
        }
    }
}

# ---------------------------------------------------------------------
# 2. New track data (No., Title, Composers, Performer, Duration)
# ---------------------------------------------------------------------
$tracks = @(
    @(1,  "Cold Little Heart",          "Brian Burton / Dean Josiah / Michael Kiwanuka", "Michael Kiwanuka", 0.4236111111111111),
    @(2,  "Black Man in a White World", "Dean Josiah / Michael Kiwanuka",                "Michael Kiwanuka", 0.17916666666666667),
    @(3,  "Falling",                    "Brian Burton / Michael Kiwanuka",               "Michael Kiwanuka", 0.17777777777777778),
    @(4,  "Place I Belong",             "Michael Kiwanuka",                              "Michael Kiwanuka", 0.19930555555555554),
    @(5,  "Love & Hate",                "Brian Burton / Dean Josiah / Michael Kiwanuka", "Michael Kiwanuka", 0.29652777777777778),
    @(6,  "One More Night",             "Michael Kiwanuka",                              "Michael Kiwanuka", 0.16180555555555556),
    @(7,  "I'll Never Love",            "Michael Kiwanuka",                              "Michael Kiwanuka", 0.11458333333333333),
    @(8,  "Rule the World",             "Brian Burton / Dean Josiah / Michael Kiwanuka", "Michael Kiwanuka", 0.23750000000000002),
    @(9,  "Father's Child",             "Brian Burton / Dean Josiah / Michael Kiwanuka", "Michael Kiwanuka", 0.2951388888888889),
    @(10, "The Final Frame",            "Brian Burton / Dean Josiah / Michael Kiwanuka", "Michael Kiwanuka", 0.2076388888888889)
)

foreach ($sheetName in @("Sheet1", "Sheet3")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $r = 2
    foreach ($t in $tracks) {
        $ws.Cells.Item($r, 1).Value = $t[0]
        $ws.Cells.Item($r, 2).Value = $t[1]
        $ws.Cells.Item($r, 3).Value = $t[2]
        $ws.Cells.Item($r, 4).Value = $t[3]
        $ws.Cells.Item($r, 5).Value = $t[4]
        $r++
    }

    # New column widths after the data refresh
    $ws.Columns.Item(1).ColumnWidth = 3.5546875
    $ws.Columns.Item(2).ColumnWidth = 25.88671875
    $ws.Columns.Item(3).ColumnWidth = 43.77734375
    $ws.Columns.Item(4).ColumnWidth = 17.109375
    $ws.Columns.Item(5).ColumnWidth = 6.109375
}

# ---------------------------------------------------------------------
# 3. Sheet2 is a fully formula-driven report off Sheet1; only its
#    selection needs to follow the data growing to row 14.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("A3:K14").Select()

$wb.Application.Calculate()
